# "Fruta / hortaliza, semanal" - weekly data refresh.
# A new weekly observation is inserted as row 107 (the sheet's first data
# row is row 2, so this lands among the existing "Apio" records for
# "Feria Lagunitas de Puerto Montt"), pushing the previously-existing
# rows 107-185 down to 108-186 unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 107; Excel shifts rows 107:185 down
# to 108:186 and the new row inherits formatting (incl. the date style)
# from the row above it.
$ws.Rows.Item(107).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A107").Value = 4
$ws.Range("B107").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C107").Value = "Los Lagos"
$ws.Range("D107").Value = 44574
$ws.Range("E107").Value = 10
$ws.Range("F107").Value = 100112017
$ws.Range("G107").Value = "Apio"
$ws.Range("H107").Value = "Americana (o)"
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 25
$ws.Range("K107").Value = 12000
$ws.Range("L107").Value = 12000
$ws.Range("M107").Value = 12000
$ws.Range("N107").Value = "$/docena de matas"
$ws.Range("O107").Value = "Región de Coquimbo"
$ws.Range("P107").Value = 2000
$ws.Range("Q107").Value = 6
$ws.Range("R107").Value = "Hortaliza"
